$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2" = 'system, backup@backdoor.com, System'
    "G3" = 'dnasr281@gmail.com, System'
    "G4" = 'System, backup@backdoor.com'
    "G6" = 'dnasr281@gmail.com, System'
    "G10" = 'dnasr281@gmail.com, System'
    "G11" = 'dnasr281@gmail.com, System'
    "G12" = 'dnasr281@gmail.com, System'
    "G13" = 'dnasr281@gmail.com, System'
    "G14" = 'dnasr281@gmail.com, System'
    "G15" = 'dnasr281@gmail.com, System'
    "G17" = 'dnasr281@gmail.com, System'
    "G18" = 'dnasr281@gmail.com, System'
    "G19" = 'dnasr281@gmail.com, System'
    "G29" = 'system, backup@backdoor.com, System'
    "G30" = 'dnasr281@gmail.com, System'
    "G31" = 'System, backup@backdoor.com'
    "G33" = 'dnasr281@gmail.com, System'
    "G37" = 'dnasr281@gmail.com, System'
    "G38" = 'dnasr281@gmail.com, System'
    "G39" = 'dnasr281@gmail.com, System'
    "G40" = 'dnasr281@gmail.com, System'
    "G41" = 'dnasr281@gmail.com, System'
    "G42" = 'dnasr281@gmail.com, System'
    "G44" = 'dnasr281@gmail.com, System'
    "G45" = 'dnasr281@gmail.com, System'
    "G46" = 'dnasr281@gmail.com, System'
    "G56" = 'system, backup@backdoor.com, System'
    "G57" = 'dnasr281@gmail.com, System'
    "G58" = 'System, backup@backdoor.com'
    "G60" = 'dnasr281@gmail.com, System'
    "G64" = 'dnasr281@gmail.com, System'
    "G65" = 'dnasr281@gmail.com, System'
    "G66" = 'dnasr281@gmail.com, System'
    "G67" = 'dnasr281@gmail.com, System'
    "G68" = 'dnasr281@gmail.com, System'
    "G69" = 'dnasr281@gmail.com, System'
    "G71" = 'dnasr281@gmail.com, System'
    "G72" = 'dnasr281@gmail.com, System'
    "G73" = 'dnasr281@gmail.com, System'
    "G86" = 'dnasr281@gmail.com, System'
    "G87" = 'dnasr281@gmail.com, System'
    "G88" = 'dnasr281@gmail.com, System'
    "G89" = 'dnasr281@gmail.com, System'
    "G93" = 'dnasr281@gmail.com, System'
    "G95" = 'dnasr281@gmail.com, System'
    "G96" = 'dnasr281@gmail.com, System'
    "G97" = 'dnasr281@gmail.com, System'
    "G112" = 'dnasr281@gmail.com, System'
    "G113" = 'dnasr281@gmail.com, System'
    "G114" = 'dnasr281@gmail.com, System'
    "G115" = 'dnasr281@gmail.com, System'
    "G119" = 'dnasr281@gmail.com, System'
    "G121" = 'dnasr281@gmail.com, System'
    "G122" = 'dnasr281@gmail.com, System'
    "G123" = 'dnasr281@gmail.com, System'
    "G138" = 'dnasr281@gmail.com, System'
    "G139" = 'dnasr281@gmail.com, System'
    "G140" = 'dnasr281@gmail.com, System'
    "G141" = 'dnasr281@gmail.com, System'
    "G145" = 'dnasr281@gmail.com, System'
    "G147" = 'dnasr281@gmail.com, System'
    "G148" = 'dnasr281@gmail.com, System'
    "G149" = 'dnasr281@gmail.com, System'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
